$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source range that gets copied: F7:I43 (LP / Uz / Ud / If table for OA95,
# rows where LP = 3 .. 39)
$src = $ws.Range("F7:I43")

# First paste: place a copy of the table at W3 (so W3:Z39 mirrors F7:I43)
$src.Copy()
$ws.Range("W3").PasteSpecial(-4163) | Out-Null

# Second paste: place another copy of the table at AB15 (so AB15:AE51
# mirrors F7:I43 as well)
$src.Copy()
$ws.Range("AB15").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0

# Update the selection to match the final state left by the paste
$ws.Range("AB15:AE51").Select() | Out-Null
